# Revert "US16 & US17"
# This reverts the added Sprint3 (sheet7) row 3 & row 4 detail data that was
# introduced by the US16/US17 commit: the extra Type/Estimate/Actual/Date/
# Assigned-To/Tested-By/Time-Spent columns (D, G:L, N:P) added to the
# "list_siblings" and "siblings_by_age" backlog-item rows are cleared out,
# restoring those rows to their original (pre-US16/US17) shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint3")

# Remove the cells (not just their values) that were added for US16 & US17
# in rows 3 and 4: column D, columns G through L, and columns N through P.
$ws.Range("D3:D4").Clear()
$ws.Range("G3:L4").Clear()
$ws.Range("N3:P4").Clear()

# The removed content made rows 3 & 4 wrap to extra lines; once the content
# is gone the rows go back to the sheet's standard (auto-fit) height.
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()

# Restore the selection/cursor position left behind on the active sheet.
$ws.Range("P9").Select() | Out-Null
